# Update "riwayat pembayaran" sheet: append 5 new payment history rows
# (rows 2-6) below the existing header row, and record the last payment
# amount (row 6, column E) as a real number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not auto-converted to a number),
# without leaving a lingering number-format style on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$rows = @(
    @{ Row = 2; NIK = "2345678990112444"; Plat = "BG6701HI"; Nama = "Tiara"; Tanggal = "01-08-2025 05:09"; Jumlah = "10000"; Metode = "BRI" },
    @{ Row = 3; NIK = "2345678990112444"; Plat = "BG6701HI"; Nama = "Tiara"; Tanggal = "01-08-2025 05:15"; Jumlah = "10000"; Metode = "BRI" },
    @{ Row = 4; NIK = "2345678990112444"; Plat = "BG6701HI"; Nama = "Tiara"; Tanggal = "01-08-2025 05:16"; Jumlah = "0";     Metode = "BRI" },
    @{ Row = 5; NIK = "2345678990112444"; Plat = "BG6701HI"; Nama = "Tiara"; Tanggal = "01-08-2025 05:24"; Jumlah = "10000"; Metode = "BRI" }
)

foreach ($r in $rows) {
    $row = $r.Row
    Set-TextValue $ws.Cells.Item($row, 1) $r.NIK
    $ws.Cells.Item($row, 2).Value = $r.Plat
    $ws.Cells.Item($row, 3).Value = $r.Nama
    $ws.Cells.Item($row, 4).Value = $r.Tanggal
    Set-TextValue $ws.Cells.Item($row, 5) $r.Jumlah
    $ws.Cells.Item($row, 6).Value = $r.Metode
}

# Row 6: the latest payment. Jumlah is recorded as a genuine number this
# time, reflecting the updated "last payment" logic.
Set-TextValue $ws.Cells.Item(6, 1) "2345678990112444"
$ws.Cells.Item(6, 2).Value = "BG6701HI"
$ws.Cells.Item(6, 3).Value = "Tiara"
$ws.Cells.Item(6, 4).Value = "01-08-2025 05:35"
$ws.Cells.Item(6, 5).Value = 10000
$ws.Cells.Item(6, 6).Value = "BRI"
